$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.587.45"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.599.11"
$ws.Range("E3").Value = "  -0.06%  "
$ws.Range("E4").Value = "  +0.25%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "592.51"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -0.24%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "154.46"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  +0.77%  "
$ws.Range("D9").Value = "2.597.70"
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  +9.98%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.159"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.62%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "5.23"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -0.35%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.352"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -2.08%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "27.41"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.59%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.0000184"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +1.54%  "
$ws.Range("D16").Value = "3.077.31"
$ws.Range("E16").Value = "  +0.98%  "
$ws.Range("D17").Value = "67.507.54"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "2.600.94"
$ws.Range("E18").Value = "  +0.31%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "11.12"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.87%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "362.65"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.98%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "7.59"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -3.57%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.30"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.93%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "2.00"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -5.23%  "
$ws.Range("E24").Value = "  -0.04%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "67.39"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "9.72"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -7.83%  "
$ws.Range("D27").Value = "2.742.05"
$ws.Range("E28").Value = "  -2.21%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "572.64"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -6.79%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "1.02"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "1.41"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -4.11%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "7.89"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -2.40%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "1.85"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.19%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.132"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -2.13%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  -5.37%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "4.92"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  -3.37%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "157.92"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +1.60%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "19.26"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.58%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.368"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.33%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.29"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "1.83"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "2.52"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -7.00%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "41.18"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.06%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.00"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "16.42"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +0.55%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "154.79"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.05%  "
$ws.Range("D48").Value = "0.0₆0285"
$ws.Range("E48").Value = "  -7.76%  "
$ws.Range("E49").Value = "  -1.43%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.626"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +1.59%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "20.64"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -4.59%  "
